$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "67.253.84"
$ws.Cells.Item(2, 5).Value = "  +4.85%  "
$ws.Cells.Item(3, 4).Value = "3.245.12"
$ws.Cells.Item(3, 5).Value = "  +2.31%  "
$ws.Cells.Item(5, 4).Value = "577.44"
$ws.Cells.Item(5, 5).Value = "  +2.34%  "
$ws.Cells.Item(6, 4).Value = "178.65"
$ws.Cells.Item(6, 5).Value = "  +5.79%  "
$ws.Cells.Item(7, 5).Value = "  -0.04%  "
$ws.Cells.Item(8, 5).Value = "  -1.30%  "
$ws.Cells.Item(9, 4).Value = "3.244.58"
$ws.Cells.Item(9, 5).Value = "  +2.28%  "
$ws.Cells.Item(10, 5).Value = "  +4.58%  "
$ws.Cells.Item(11, 5).Value = "  +2.11%  "
$ws.Cells.Item(12, 5).Value = "  +4.45%  "
$ws.Cells.Item(13, 4).Value = "3.806.93"
$ws.Cells.Item(13, 5).Value = "  +2.20%  "
$ws.Cells.Item(14, 5).Value = "  +0.74%  "
$ws.Cells.Item(15, 4).Value = "27.94"
$ws.Cells.Item(15, 5).Value = "  +2.24%  "
$ws.Cells.Item(16, 4).Value = "67.181.64"
$ws.Cells.Item(16, 5).Value = "  +4.73%  "
$ws.Cells.Item(17, 5).Value = "  +3.20%  "
$ws.Cells.Item(18, 4).Value = "3.245.83"
$ws.Cells.Item(18, 5).Value = "  +2.22%  "
$ws.Cells.Item(19, 5).Value = "  +2.35%  "
$ws.Cells.Item(20, 5).Value = "  +2.69%  "
$ws.Cells.Item(21, 4).Value = "373.65"
$ws.Cells.Item(21, 5).Value = "  +6.05%  "
$ws.Cells.Item(22, 5).Value = "  +6.00%  "
$ws.Cells.Item(23, 5).Value = "  +0.11%  "
$ws.Cells.Item(24, 4).Value = "71.12"
$ws.Cells.Item(24, 5).Value = "  +3.77%  "
$ws.Cells.Item(25, 5).Value = "  +0.91%  "
$ws.Cells.Item(26, 4).Value = "3.382.32"
$ws.Cells.Item(26, 5).Value = "  +2.13%  "
$ws.Cells.Item(27, 4).Value = "0.0000119"
$ws.Cells.Item(27, 5).Value = "  +0.60%  "
$ws.Cells.Item(28, 4).Value = "9.85"
$ws.Cells.Item(28, 5).Value = "  +2.88%  "
$ws.Cells.Item(29, 5).Value = "  +2.26%  "
$ws.Cells.Item(30, 5).Value = "  +0.37%  "
$ws.Cells.Item(31, 4).Value = "1.98"
$ws.Cells.Item(31, 5).Value = "  +4.13%  "
$ws.Cells.Item(33, 4).Value = "22.56"
$ws.Cells.Item(33, 5).Value = "  +2.86%  "
$ws.Cells.Item(34, 2).Value = "Fetch.AI"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(34, 4).Value = "1.28"
$ws.Cells.Item(34, 5).Value = "  +6.82%  "
$ws.Cells.Item(35, 2).Value = "USDe"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(35, 4).Value = "0.998"
$ws.Cells.Item(35, 5).Value = "  +0.03%  "
$ws.Cells.Item(36, 4).Value = "6.82"
$ws.Cells.Item(36, 5).Value = "  +2.82%  "
$ws.Cells.Item(37, 2).Value = "Monero"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(37, 4).Value = "162.45"
$ws.Cells.Item(37, 5).Value = "  +5.34%  "
$ws.Cells.Item(38, 2).Value = "ImmutableX"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(38, 4).Value = "1.50"
$ws.Cells.Item(38, 5).Value = "  +4.43%  "
$ws.Cells.Item(39, 5).Value = "  +5.12%  "
$ws.Cells.Item(40, 4).Value = "1.86"
$ws.Cells.Item(40, 5).Value = "  +10.16%  "
$ws.Cells.Item(41, 4).Value = "6.88"
$ws.Cells.Item(41, 5).Value = "  +14.88%  "
$ws.Cells.Item(42, 4).Value = "26.82"
$ws.Cells.Item(42, 5).Value = "  +3.17%  "
$ws.Cells.Item(43, 4).Value = "2.62"
$ws.Cells.Item(43, 5).Value = "  +5.66%  "
$ws.Cells.Item(44, 4).Value = "2.764.96"
$ws.Cells.Item(44, 5).Value = "  +5.85%  "
$ws.Cells.Item(45, 4).Value = "357.61"
$ws.Cells.Item(45, 5).Value = "  +12.76%  "
$ws.Cells.Item(46, 5).Value = "  +5.60%  "
$ws.Cells.Item(47, 4).Value = "25.83"
$ws.Cells.Item(47, 5).Value = "  +9.49%  "
$ws.Cells.Item(48, 4).Value = "40.45"
$ws.Cells.Item(48, 5).Value = "  +3.09%  "
$ws.Cells.Item(49, 4).Value = "0.0675"
$ws.Cells.Item(49, 5).Value = "  +4.73%  "
$ws.Cells.Item(50, 5).Value = "  +3.33%  "
$ws.Cells.Item(51, 5).Value = "  +1.50%  "
